# Adds the new match row (row 79) to the sheet, replicating the data
# produced by the scraping script run on 01-12-2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

$ws.Cells.Item($row, 1).Value = 78
$ws.Cells.Item($row, 2).Value = "croatia"
$ws.Cells.Item($row, 3).Value = "hnl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45261.70833333334
$ws.Cells.Item($row, 6).Value = "Istra 1961"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Varazdin"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 2.31
$ws.Cells.Item($row, 11).Value = "26/11/2023 15:12"
$ws.Cells.Item($row, 12).Value = 2.21
$ws.Cells.Item($row, 13).Value = "01/12/2023 16:52"
$ws.Cells.Item($row, 14).Value = 3.33
$ws.Cells.Item($row, 15).Value = "26/11/2023 15:12"
$ws.Cells.Item($row, 16).Value = 3.3
$ws.Cells.Item($row, 17).Value = "01/12/2023 16:50"
$ws.Cells.Item($row, 18).Value = 3.13
$ws.Cells.Item($row, 19).Value = "26/11/2023 15:12"
$ws.Cells.Item($row, 20).Value = 3.29
$ws.Cells.Item($row, 21).Value = "01/12/2023 16:45"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/istra-1961-varazdin/fB7HUU9T/"

# Match the style of column A (bold/bordered index style) and column E (datetime format)
# by copying the formatting from the row above (xlPasteFormats = -4122)
$ws.Range("A" + ($row - 1)).Copy() | Out-Null
$ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null

$ws.Range("E" + ($row - 1)).Copy() | Out-Null
$ws.Range("E" + $row).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
